$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 91 (shifts rows 91-109 down to 92-110),
# copying formatting from the row above.
$ws.Rows(91).Insert()

# Populate the newly inserted row with the new parameter.
$ws.Range("A91").Value = "success_rate_surgical_removal_placenta"
$ws.Range("B91").Value = 0.7

# Column B width nudges slightly wider in the target workbook.
$ws.Columns(2).ColumnWidth = 25.5

# Update the selection / active cell to match the saved view state.
$ws.Range("D92").Select()
